$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 58: Ikbel Hadj Hassine ---
$ws.Range("A58").Value = "Ikbel"
$ws.Range("B58").Value = "Hadj Hassine"
$ws.Range("C58").Value = "Université de Monastir"
$ws.Range("D58").Value = "Tunisie"
$ws.Range("E58").Value = "iGnXpRMAAAAJ"
$ws.Range("F58").Value = "F"
$ws.Range("G58").Value = 1990
$ws.Range("H58").Value = "Médecine, Biologie et Sciences de la Santé"

# --- Row 59: Raouia Mokni ---
$ws.Range("A59").Value = "Raouia"
$ws.Range("B59").Value = "Mokni"
$ws.Range("C59").Value = "Université de Gabès"
$ws.Range("D59").Value = "Tunisie"
$ws.Range("E59").Value = "j_a72EQAAAAJ"
$ws.Range("F59").Value = "F"
$ws.Range("G59").Value = 1986
$ws.Range("H59").Value = "Informatique, Mathématiques et Ingénierie"

# --- Row 60: Nesrine Zitouni ---
$ws.Range("A60").Value = "Nesrine"
$ws.Range("B60").Value = "Zitouni"
$ws.Range("C60").Value = "Université de Caen Normandie"
$ws.Range("D60").Value = "France"
$ws.Range("E60").Value = "kBOIwKQAAAAJ"
$ws.Range("F60").Value = "F"
$ws.Range("G60").Value = 1991
$ws.Range("H60").Value = "Médecine, Biologie et Sciences de la Santé"

# --- Fix row 55 (Mohamed Raâfet Ben Khedher): institution / country ---
$ws.Range("C55").Value = "Université de Jendouba"
$ws.Range("D55").Value = "Tunisie"

# --- Row 61: Hajer Bougatef ---
$ws.Range("A61").Value = "Hajer"
$ws.Range("B61").Value = "Bougatef"
$ws.Range("C61").Value = "Université de Sfax"
$ws.Range("D61").Value = "Tunisie"
$ws.Range("E61").Value = "ocEvNeAAAAAJ"
$ws.Range("F61").Value = "F"
$ws.Range("G61").Value = 1993
$ws.Range("H61").Value = "Médecine, Biologie et Sciences de la Santé"

# --- Row 62: Melek Hajji ---
$ws.Range("A62").Value = "Melek"
$ws.Range("B62").Value = "Hajji"
$ws.Range("C62").Value = "Université de Kairouan"
$ws.Range("D62").Value = "Tunisie"
$ws.Range("E62").Value = "DJHqHkgAAAAJ"
$ws.Range("F62").Value = "M"
$ws.Range("G62").Value = 1989
$ws.Range("H62").Value = "Chimie et Sciences des Matériaux"

# --- Fix row 57 (Latifa Remadi): institution / country ---
$ws.Range("C57").Value = "Université de Monastir"
$ws.Range("D57").Value = "Tunisie"

# Apply the same "Genre" column style (Arial 8pt grey) used throughout the
# sheet to the newly-added cells, by copying the format from F2.
$ws.Range("F2").Copy()
$ws.Range("F58:F62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view/selection state ---
$ws.Range("F63").Select()
